$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.975.38"
$ws.Range("E2").Value = "  +5.43%  "
$ws.Range("D3").Value = "1.915.14"
$ws.Range("E3").Value = "  +4.65%  "
$ws.Range("E4").Value = "  -0.28%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "338.95"
$cell.Style = "Normal"
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "1.000"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -0.26%  "
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.4745"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  +3.58%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.4064"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  +6.77%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "48.17"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +3.92%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.08174"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +3.34%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "1.029"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +6.03%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "22.49"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +6.69%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.899.30"
$ws.Range("E13").Value = "  +4.68%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "6.094"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +3.62%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "7.385"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +4.40%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "91.61"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +2.32%  "
$ws.Range("E17").Value = "  -0.28%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "0.00001051"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +2.37%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "0.06631"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -0.09%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "17.85"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +4.27%  "
$ws.Range("E21").Value = "  -0.47%  "
$ws.Range("D22").Value = "29.001.47"
$ws.Range("E22").Value = "  +5.60%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "5.583"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +4.60%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "11.21"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +3.57%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "2.268"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -1.04%  "
$ws.Range("D26").Value = "2.128.30"
$ws.Range("E26").Value = "  +4.79%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "160.76"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +3.29%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "20.05"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +3.34%  "
$ws.Range("E29").Value = "  +5.68%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "5.532"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +4.41%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "121.09"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +2.18%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "1.011"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +7.32%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "0.09595"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +3.15%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "1.435"
$cell.Style = "Normal"
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "3.641"
$cell.Style = "Normal"
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "5.433"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +3.46%  "
$ws.Range("E37").Value = "  +4.55%  "
$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "8.689"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +7.56%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.02277"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +4.45%  "
$ws.Range("E40").Value = "  +5.50%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.6040"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +4.60%  "
$ws.Range("E42").Value = "  +6.02%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.1901"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +3.97%  "
$ws.Range("E44").Value = "  -0.24%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "12.41"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +3.94%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.5637"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +3.41%  "
$ws.Range("E48").Value = "  +5.93%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "0.07252"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +9.68%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "2.162"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +19.77%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "113.54"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +2.24%  "
